$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds a "Date" column that was incorrectly populated with a
# label derived from the source file name ("6-25-2007-08") for every data
# row. The NBA stats for this file actually correspond to the date
# 2008-06-25 (the stats were captured one calendar day off because of how
# the league reports/labels games), so replace the mislabeled values with
# the correct ISO-style date string.
$oldValue = "6-25-2007-08"
$newValue = "2008-06-25"
$dateColumn = 58  # column BF

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateColumn)
    if ($cell.Value() -eq $oldValue) {
        # Force the cell to remain plain text so Excel does not reinterpret
        # "2008-06-25" as a date serial value, then strip the temporary
        # text number-format again so no extra cell style is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    }
}
